## Auto-committed on 2023/06/02 週五 17:31:38.49
## Adds a new "AcDate" field row to the DBD field list, renumbers the
## following SEQ values, and records the new findAcDate read-key example
## on the DBS sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "DBD" sheet - insert a new field row (AcDate / 會計日期 / DecimalD)
#    right above the existing "CompanyName" row (row 9).
# ---------------------------------------------------------------------
$dbd = $wb.Worksheets.Item("DBD")

$dbd.Rows(9).Insert()

$dbd.Range("A9").Value = 1
$dbd.Range("B9").Value = "AcDate"
$dbd.Range("C9").Value = "會計日期"
$dbd.Range("D9").Value = "DecimalD"
$dbd.Range("E9").Value = 8

# Renumber the SEQ column (A) for every row that got pushed down one
# slot so the numbering stays contiguous (2..11).
$dbd.Range("A10").Value = 2
$dbd.Range("A11").Value = 3
$dbd.Range("A12").Value = 4
$dbd.Range("A13").Value = 5
$dbd.Range("A14").Value = 6
$dbd.Range("A15").Value = 7
$dbd.Range("A16").Value = 8
$dbd.Range("A17").Value = 9
$dbd.Range("A18").Value = 10
$dbd.Range("A19").Value = 11

$dbd.Range("B12").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) "DBS" sheet - document the new findAcDate read-key example.
# ---------------------------------------------------------------------
$dbs = $wb.Worksheets.Item("DBS")

$dbs.Range("A2").Value = "findAcDate"
$dbs.Range("B2").Value = "AcDate = "

# ---------------------------------------------------------------------
# 3) "SP" sheet - no longer the active tab after this session.
# ---------------------------------------------------------------------
$dbd.Select() | Out-Null
